$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.882.29"
$ws.Range("D3").Value = "1.806.22"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'310.27"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.4445"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("D8").Value = "'0.3674"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.07327"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'0.8563"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'20.65"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.861.02"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "'92.43"
$ws.Range("E14").Value = "  +3.25%  "
$ws.Range("D15").Value = "'5.304"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "'0.07064"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "26.924.05"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "'5.150"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'10.83"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'1.993"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'151.93"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.187"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.47"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'5.215"
$ws.Range("D29").Value = "'116.65"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "'0.08838"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'0.7499"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "'2.935"
$ws.Range("E33").Value = "  +4.63%  "
$ws.Range("D34").Value = "'4.462"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "'0.9994"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'0.01969"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'0.05201"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "'0.5321"
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("D40").Value = "'2.861"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "'7.016"
$ws.Range("E41").Value = "  -3.93%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.5165"
$ws.Range("E43").Value = "  +8.99%  "
$ws.Range("D44").Value = "'8.423"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "'1.984"
$ws.Range("E45").Value = "  +6.88%  "
$ws.Range("D46").Value = "'10.53"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "'105.34"
$ws.Range("D48").Value = "'0.9988"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'1.668"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'0.06320"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "'0.9179"
$ws.Range("E51").Value = "  +0.52%  "
